$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the text "R40". The rule row now needs the literal
# text "1" instead, while keeping the cell's existing General-formatted
# style (s="23") untouched. A direct Value assignment of "1" would be
# auto-coerced to a number by Excel's type inference, so we instead write
# a formula that evaluates to the text string "1" and then convert the
# formula result to a static value via copy / paste-special (values only).
# This preserves the cell's number format / style while still producing a
# literal text cell (stored as a shared string) rather than a number.
$ws.Range("B11").Formula = "=""1"""
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)
$excel.CutCopyMode = $false
